# "remove useless fields of template"
#
# The template used helper cells M9:O9 (labels "Толщина"/"Ширина"/"Цвет"),
# M10:O10 (VLOOKUP formulas pulling thickness/width/color from the external
# "Данные" workbook) and M11 (a VLOOKUP formula) plus P11 (the literal
# "белый" color swatch) purely as scratch/helper fields. They are not shown
# anywhere else on the sheet except being fed into the C11 CONCATENATE
# formula (through M10/N10). These helper fields are unused by the final
# document, so clear them out entirely (values + formulas), leaving the
# cells blank but keeping their existing formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 helper header labels ("Толщина", "Ширина", "Цвет")
$ws.Range("M9:P9").ClearContents()

# Row 10 helper VLOOKUP formulas (thickness/width/color lookups)
$ws.Range("M10:P10").ClearContents()

# Row 11 helper VLOOKUP formula and the leftover "белый" literal
$ws.Range("M11:P11").ClearContents()

$wb.Save()
